$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts existing D:K data to F:M),
# matching the two additional quarterly periods added to each of the three tables.
$ws.Range("D1:E1").EntireColumn.Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)

# Propagate the number formatting from the (now-shifted) old D:E columns, which now
# live at F:G, into the freshly inserted D:E columns so the new cells keep the same look
# (date format row, thousands style, etc.) as the rest of the table.
$ws.Range("F5:G102").Copy() | Out-Null
$ws.Range("D5:E102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Cells that should read "NA" in the new columns (matches the existing NA cells in this row)
$naCells = "D12,E12,D29,E29"
$ws.Range($naCells).Value = "NA"

# Cells that should be zero in the new columns (matches the existing 0 cells in this row)
$zeroCells = "D13,E13,D14,E14,D25,E25,D28,E28,D30,E30,D31,E31,D34,E34,D42,E42,D50,E50,D51,E51,D53,E53,D63,E63,D64,E64,D65,E65,D68,E68,D69,E69,D70,E70,D71,E71,D72,E72,D73,E73,D74,E74,D75,E75,D77,E77,D84,E84,D85,E85,D86,E86,D87,E87,D88,E88,D92,E92,D93,E93,D97,E97,D98,E98,D99,E99,E101"
$ws.Range($zeroCells).Value = 0

# New quarterly figures for the two newly added periods (columns D and E)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 8786000
$ws.Range("E8").Value = 8792000
$ws.Range("D9").Value = 6955000
$ws.Range("E9").Value = 7768000
$ws.Range("D10").Value = 1831000
$ws.Range("E10").Value = 1024000
$ws.Range("D15").Value = 214000
$ws.Range("E15").Value = 131000
$ws.Range("D17").Value = 7472000
$ws.Range("E17").Value = 8300000
$ws.Range("D18").Value = 1314000
$ws.Range("E18").Value = 492000
$ws.Range("D20").Value = 67000
$ws.Range("E20").Value = 315000
$ws.Range("D21").Value = 1595000
$ws.Range("E21").Value = 938000
$ws.Range("D22").Value = 102000
$ws.Range("E22").Value = 108000
$ws.Range("D23").Value = 1279000
$ws.Range("E23").Value = 699000
$ws.Range("D24").Value = 218000
$ws.Range("E24").Value = 23000
$ws.Range("D26").Value = 1061000
$ws.Range("E26").Value = 676000
$ws.Range("D27").Value = 180000
$ws.Range("E27").Value = 111000
$ws.Range("D32").Value = -67000
$ws.Range("E32").Value = -315000
$ws.Range("D33").Value = 180000
$ws.Range("E33").Value = 111000
$ws.Range("D35").Value = 180000
$ws.Range("E35").Value = 111000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 69000
$ws.Range("E41").Value = 33000
$ws.Range("D43").Value = 2454000
$ws.Range("E43").Value = 2954000
$ws.Range("D44").Value = 640000
$ws.Range("E44").Value = 824000
$ws.Range("D45").Value = 373000
$ws.Range("E45").Value = 320000
$ws.Range("D46").Value = 3536000
$ws.Range("E46").Value = 4131000
$ws.Range("D47").Value = 2702000
$ws.Range("E47").Value = 2539000
$ws.Range("D48").Value = 15718000
$ws.Range("E48").Value = 15606000
$ws.Range("D49").Value = 2521000
$ws.Range("E49").Value = 2540000
$ws.Range("D52").Value = 2353000
$ws.Range("E52").Value = 2474000
$ws.Range("D54").Value = 26830000
$ws.Range("E54").Value = 27290000
$ws.Range("D57").Value = 2705000
$ws.Range("E57").Value = 3614000
$ws.Range("D58").Value = 66000
$ws.Range("E58").Value = 429000
$ws.Range("D59").Value = 687000
$ws.Range("E59").Value = 615000
$ws.Range("D60").Value = 3458000
$ws.Range("E60").Value = 4658000
$ws.Range("D61").Value = 9143000
$ws.Range("E61").Value = 9140000
$ws.Range("D62").Value = 910000
$ws.Range("E62").Value = 781000
$ws.Range("D66").Value = 24984000
$ws.Range("E66").Value = 25543000
$ws.Range("D76").Value = 1846000
$ws.Range("E76").Value = 1747000
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 180000
$ws.Range("E81").Value = 111000
$ws.Range("D83").Value = 214000
$ws.Range("E83").Value = 131000
$ws.Range("D89").Value = 1312000
$ws.Range("E89").Value = 279000
$ws.Range("D91").Value = 1139000
$ws.Range("E91").Value = -460000
$ws.Range("D94").Value = -629000
$ws.Range("E94").Value = 322000
$ws.Range("D96").Value = -47000
$ws.Range("E96").Value = -48000
$ws.Range("D100").Value = -641000
$ws.Range("E100").Value = -605000
$ws.Range("D101").Value = -6000
$ws.Range("D102").Value = 36000
$ws.Range("E102").Value = -4000

# Data correction carried in with the new quarter: "Net Income" row (H91) was restated
$ws.Range("H91").Value = 776000

